$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "G4"
$ws.Range("B5").Value = "Read Book"
$ws.Range("C5").Value = "Daily"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 45907
$ws.Range("E5").NumberFormat = "YYYY-MM-DD"
$ws.Range("F5").Value = 36
